$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, D(serial date), L(Calidad), M(Volumen), N(Precio minimo), O(Precio maximo),
# P(Precio promedio ponderado), Q(Unidad de comercializacion), R(Origen), S(Precio $/Kg), T(Kg/unidad)
# Values correspond to the original row noted in the comment (rows were re-shuffled).
$rowsData = @(
    ,@(2, 44270, "Primera", 85, 12000, 12000, 12000, "$/caja 14 kilos granel", "Provincia del Elquí", 857, 14)  # from row 41
    ,@(3, 44252, "Primera", 60, 14000, 14000, 14000, "$/caja 14 kilos empedrada", "Provincia de Limarí", 1000, 14)  # from row 11
    ,@(4, 45001, "Primera", 50, 16000, 16000, 16000, "$/caja 14 kilos empedrada", "Provincia de Limarí", 1143, 14)  # from row 28
    ,@(5, 45015, "Primera", 56, 15000, 15000, 15000, "$/caja 14 kilos empedrada", "Provincia de Limarí", 1071, 14)  # from row 17
    ,@(6, 45042, "Especial", 50, 17000, 17000, 17000, "$/caja 14 kilos granel", "Provincia de Limarí", 1214, 14)  # from row 24
    ,@(7, 45042, "Primera", 50, 14000, 14000, 14000, "$/caja 14 kilos granel", "Provincia de Limarí", 1000, 14)  # from row 25
    ,@(8, 44614, "Primera", 54, 14000, 14000, 14000, "$/caja 14 kilos granel", "Provincia de Limarí", 1000, 14)  # from row 6
    ,@(9, 44271, "Primera", 50, 12000, 12000, 12000, "$/caja 14 kilos granel", "Provincia del Elquí", 857, 14)  # from row 3
    ,@(10, 44627, "Primera", 56, 17000, 17000, 17000, "$/caja 14 kilos empedrada", "Provincia de Limarí", 1214, 14)  # from row 12
    ,@(11, 44323, "Primera", 60, 14000, 14000, 14000, "$/caja 14 kilos granel", "Provincia de Limarí", 1000, 14)  # from row 33
    ,@(12, 44239, "Primera", 70, 15000, 15000, 15000, "$/caja 15 kilos granel", "Provincia de Limarí", 1000, 15)  # from row 29
    ,@(13, 44312, "Primera", 68, 14000, 14000, 14000, "$/caja 14 kilos granel", "Provincia de Limarí", 1000, 14)  # from row 13
    ,@(14, 45054, "Especial", 54, 16000, 16000, 16000, "$/caja 14 kilos empedrada", "Provincia de Limarí", 1143, 14)  # from row 4
    ,@(15, 45054, "Primera", 50, 14000, 14000, 14000, "$/caja 14 kilos empedrada", "Provincia de Limarí", 1000, 14)  # from row 5
    ,@(16, 44260, "Primera", 56, 13000, 13000, 13000, "$/caja 14 kilos empedrada", "Provincia del Elquí", 929, 14)  # from row 9
    ,@(17, 44588, "Primera", 85, 19000, 20000, 19529, "$/caja 14 kilos granel", "Provincia de Limarí", 1395, 14)  # from row 40
    ,@(18, 44259, "Primera", 80, 12000, 12000, 12000, "$/caja 15 kilos empedrada", "Provincia de Limarí", 800, 15)  # from row 2
    ,@(19, 44315, "Primera", 65, 14000, 14000, 14000, "$/caja 14 kilos granel", "Provincia de Limarí", 1000, 14)  # from row 35
    ,@(20, 44320, "Primera", 45, 14000, 14000, 14000, "$/caja 14 kilos granel", "Provincia de Limarí", 1000, 14)  # from row 31
    ,@(21, 44313, "Primera", 36, 14000, 14000, 14000, "$/caja 14 kilos granel", "Provincia de Limarí", 1000, 14)  # from row 27
    ,@(22, 45043, "Especial", 45, 17000, 17000, 17000, "$/caja 14 kilos granel", "Provincia de Limarí", 1214, 14)  # from row 38
    ,@(23, 45043, "Primera", 67, 14000, 14000, 14000, "$/caja 14 kilos granel", "Provincia de Limarí", 1000, 14)  # from row 39
    ,@(24, 45050, "Especial", 56, 14000, 14000, 14000, "$/caja 14 kilos granel", "Provincia de Limarí", 1000, 14)  # from row 19
    ,@(25, 45050, "Primera", 50, 12000, 12000, 12000, "$/caja 14 kilos granel", "Provincia de Limarí", 857, 14)  # from row 20
    ,@(26, 44322, "Primera", 50, 14000, 14000, 14000, "$/caja 14 kilos granel", "Provincia de Limarí", 1000, 14)  # from row 21
    ,@(27, 44278, "Primera", 45, 13000, 13000, 13000, "$/caja 14 kilos empedrada", "Provincia del Elquí", 929, 14)  # from row 10
    ,@(28, 44314, "Primera", 56, 14000, 14000, 14000, "$/caja 14 kilos granel", "Provincia de Limarí", 1000, 14)  # from row 15
    ,@(29, 45044, "Especial", 30, 16000, 16000, 16000, "$/caja 14 kilos granel", "Provincia de Limarí", 1143, 14)  # from row 36
    ,@(30, 45044, "Primera", 30, 14000, 14000, 14000, "$/caja 14 kilos granel", "Provincia de Limarí", 1000, 14)  # from row 37
    ,@(31, 44316, "Primera", 48, 14000, 14000, 14000, "$/caja 14 kilos granel", "Provincia de Limarí", 1000, 14)  # from row 22
    ,@(32, 44592, "Primera", 54, 20000, 20000, 20000, "$/caja 15 kilos empedrada", "Provincia de Limarí", 1333, 15)  # from row 23
    ,@(33, 44630, "Primera", 75, 15000, 15000, 15000, "$/caja 14 kilos empedrada", "Provincia de Limarí", 1071, 14)  # from row 26
    ,@(34, 45014, "Primera", 60, 15000, 15000, 15000, "$/caja 14 kilos empedrada", "Provincia de Limarí", 1071, 14)  # from row 18
    ,@(35, 44616, "Primera", 70, 14000, 14000, 14000, "$/caja 14 kilos empedrada", "Provincia de Limarí", 1000, 14)  # from row 16
    ,@(36, 45040, "Especial", 65, 17000, 17000, 17000, "$/caja 14 kilos granel", "Provincia de Limarí", 1214, 14)  # from row 7
    ,@(37, 45040, "Primera", 60, 14000, 14000, 14000, "$/caja 14 kilos granel", "Provincia de Limarí", 1000, 14)  # from row 8
    ,@(38, 44242, "Primera", 45, 12000, 12000, 12000, "$/caja 15 kilos granel", "Provincia de Limarí", 800, 15)  # from row 34
    ,@(39, 45006, "Primera", 40, 16000, 16000, 16000, "$/caja 14 kilos empedrada", "Provincia de Limarí", 1143, 14)  # from row 14
    ,@(40, 44238, "Primera", 60, 15000, 15000, 15000, "$/caja 15 kilos granel", "Provincia de Limarí", 1000, 15)  # from row 32
    ,@(41, 44245, "Primera", 50, 15000, 15000, 15000, "$/caja 15 kilos granel", "Provincia de Limarí", 1000, 15)  # from row 30
)

foreach ($entry in $rowsData) {
    $r = $entry[0]
    $ws.Range("D$r").Value = $entry[1]
    $ws.Range("L$r").Value = $entry[2]
    $ws.Range("M$r").Value = $entry[3]
    $ws.Range("N$r").Value = $entry[4]
    $ws.Range("O$r").Value = $entry[5]
    $ws.Range("P$r").Value = $entry[6]
    $ws.Range("Q$r").Value = $entry[7]
    $ws.Range("R$r").Value = $entry[8]
    $ws.Range("S$r").Value = $entry[9]
    $ws.Range("T$r").Value = $entry[10]
}

